# Collapse the 'Persona' block (rows 4-8) and shift the 'Denuncia' block
# (old rows 9-25) up into rows 4-20, then drop the now-empty trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: denuncia_id
$ws.Cells.Item(4, 1).Value = 'Denuncias'
$ws.Cells.Item(4, 2).Value = 'Denuncia'
$ws.Cells.Item(4, 3).Value = 'denuncia_id'
$ws.Cells.Item(4, 4).Value = 'alfanumerico'
$ws.Cells.Item(4, 5).Value = 'Identificador de la denuncia realizada.'
$ws.Cells.Item(4, 6).Value = "'" + '3384556'
$ws.Cells.Item(4, 6).Style = "Normal"
$ws.Cells.Item(4, 7).Value = 'schema:Thing'

# Row 5: denuncia_fecha
$ws.Cells.Item(5, 1).Value = 'Denuncias'
$ws.Cells.Item(5, 2).Value = 'Denuncia'
$ws.Cells.Item(5, 3).Value = 'denuncia_fecha'
$ws.Cells.Item(5, 4).Value = 'fecha'
$ws.Cells.Item(5, 5).Value = 'Fecha de la denuncia.'
$ws.Cells.Item(5, 6).Value = "'" + '2019-03-12'
$ws.Cells.Item(5, 6).Style = "Normal"
$ws.Cells.Item(5, 7).Value = 'schema:Thing'

# Row 6: denuncia_hora
$ws.Cells.Item(6, 1).Value = 'Denuncias'
$ws.Cells.Item(6, 2).Value = 'Denuncia'
$ws.Cells.Item(6, 3).Value = 'denuncia_hora'
$ws.Cells.Item(6, 4).Value = 'fecha'
$ws.Cells.Item(6, 5).Value = 'Hora de la denuncia.'
$ws.Cells.Item(6, 6).Value = '14:45'
$ws.Cells.Item(6, 7).Value = 'schema:Thing'

# Row 7: denuncia_medio
$ws.Cells.Item(7, 1).Value = 'Denuncias'
$ws.Cells.Item(7, 2).Value = 'Denuncia'
$ws.Cells.Item(7, 3).Value = 'denuncia_medio'
$ws.Cells.Item(7, 4).Value = 'alfanumerico'
$ws.Cells.Item(7, 5).Value = 'Indica de que forma fue realizada la denuncia'
$ws.Cells.Item(7, 6).Value = 'Página web'
$ws.Cells.Item(7, 7).Value = ""

# Row 8: denuncia_lugar_radicacion
$ws.Cells.Item(8, 1).Value = 'Denuncias'
$ws.Cells.Item(8, 2).Value = 'Denuncia'
$ws.Cells.Item(8, 3).Value = 'denuncia_lugar_radicacion'
$ws.Cells.Item(8, 4).Value = 'alfanumerico'
$ws.Cells.Item(8, 5).Value = 'Oficina o dependencia que recepciono la denuncia ya sea presencial, electrónica, telefónica u otro medio.'
$ws.Cells.Item(8, 6).Value = ""
$ws.Cells.Item(8, 7).Value = 'schema:Thing'

# Row 9: denuncia_enlace_seguimiento
$ws.Cells.Item(9, 1).Value = 'Denuncias'
$ws.Cells.Item(9, 2).Value = 'Denuncia'
$ws.Cells.Item(9, 3).Value = 'denuncia_enlace_seguimiento'
$ws.Cells.Item(9, 4).Value = 'url'
$ws.Cells.Item(9, 5).Value = 'De existir un enlace para hacer seguimiento a la denuncia, opcional'
$ws.Cells.Item(9, 6).Value = 'https://www.bahia.gob.ar/vecinos/'
$ws.Cells.Item(9, 7).Value = 'schema:Thing'

# Row 10: denuncia_estado_actual
$ws.Cells.Item(10, 1).Value = 'Denuncias'
$ws.Cells.Item(10, 2).Value = 'Denuncia'
$ws.Cells.Item(10, 3).Value = 'denuncia_estado_actual'
$ws.Cells.Item(10, 4).Value = 'alfanumerico'
$ws.Cells.Item(10, 5).Value = 'Ultimo estado de la denuncia'
$ws.Cells.Item(10, 6).Value = 'Activo'
$ws.Cells.Item(10, 7).Value = ""

# Row 11: denuncia_direccion
$ws.Cells.Item(11, 1).Value = 'Denuncias'
$ws.Cells.Item(11, 2).Value = 'Denuncia'
$ws.Cells.Item(11, 3).Value = 'denuncia_direccion'
$ws.Cells.Item(11, 4).Value = 'alfanumerico'
$ws.Cells.Item(11, 5).Value = 'Dirección donde sucedió el hecho denunciado.'
$ws.Cells.Item(11, 6).Value = 'Alsina 1600'
$ws.Cells.Item(11, 7).Value = 'schema:PostalAddress'

# Row 12: denuncia_latitud
$ws.Cells.Item(12, 1).Value = 'Denuncias'
$ws.Cells.Item(12, 2).Value = 'Denuncia'
$ws.Cells.Item(12, 3).Value = 'denuncia_latitud'
$ws.Cells.Item(12, 4).Value = 'numerico'
$ws.Cells.Item(12, 5).Value = 'Latitud donde sucedió el hecho denunciado.'
$ws.Cells.Item(12, 6).Value = "'" + '-38.705048'
$ws.Cells.Item(12, 6).Style = "Normal"
$ws.Cells.Item(12, 7).Value = 'schema:GeoCoordinates'

# Row 13: denuncia_longitud
$ws.Cells.Item(13, 1).Value = 'Denuncias'
$ws.Cells.Item(13, 2).Value = 'Denuncia'
$ws.Cells.Item(13, 3).Value = 'denuncia_longitud'
$ws.Cells.Item(13, 4).Value = 'numerico'
$ws.Cells.Item(13, 5).Value = 'Longitud donde sucedió el hecho denunciado.'
$ws.Cells.Item(13, 6).Value = "'" + '-62.250596'
$ws.Cells.Item(13, 6).Style = "Normal"
$ws.Cells.Item(13, 7).Value = 'schema:GeoCoordinates'

# Row 14: denuncia_detalle
$ws.Cells.Item(14, 1).Value = 'Denuncias'
$ws.Cells.Item(14, 2).Value = 'Denuncia'
$ws.Cells.Item(14, 3).Value = 'denuncia_detalle'
$ws.Cells.Item(14, 4).Value = 'texto'
$ws.Cells.Item(14, 5).Value = 'Detalle brindado por el denunciante al realizar la denuncia.'
$ws.Cells.Item(14, 6).Value = 'Cable se encuentra colgando desde la columna de alumbrado.'
$ws.Cells.Item(14, 7).Value = 'schema:Thing'

# Row 15: denuncia_categoria
$ws.Cells.Item(15, 1).Value = 'Denuncias'
$ws.Cells.Item(15, 2).Value = 'Denuncia'
$ws.Cells.Item(15, 3).Value = 'denuncia_categoria'
$ws.Cells.Item(15, 4).Value = 'alfanumerico'
$ws.Cells.Item(15, 5).Value = 'Categoría o clasificación del hecho denunciado.'
$ws.Cells.Item(15, 6).Value = 'Alumbrado'
$ws.Cells.Item(15, 7).Value = 'schema:Thing'

# Row 16: denuncia_tipo
$ws.Cells.Item(16, 1).Value = 'Denuncias'
$ws.Cells.Item(16, 2).Value = 'Denuncia'
$ws.Cells.Item(16, 3).Value = 'denuncia_tipo'
$ws.Cells.Item(16, 4).Value = 'alfanumerico'
$ws.Cells.Item(16, 5).Value = 'Tipo de denuncia (esto es una clasificación más específica que "categoría").'
$ws.Cells.Item(16, 6).Value = 'Cable Suelto'
$ws.Cells.Item(16, 7).Value = 'schema:Thing'

# Row 17: denuncia_subtipo
$ws.Cells.Item(17, 1).Value = 'Denuncias'
$ws.Cells.Item(17, 2).Value = 'Denuncia'
$ws.Cells.Item(17, 3).Value = 'denuncia_subtipo'
$ws.Cells.Item(17, 4).Value = 'alfanumerico'
$ws.Cells.Item(17, 5).Value = 'Subtipo de denuncia (esto es una clasificación más específica que "tipo", si aplica).'
$ws.Cells.Item(17, 6).Value = ""
$ws.Cells.Item(17, 7).Value = 'schema:Thing'

# Row 18: denuncia_otra_clasificacion
$ws.Cells.Item(18, 1).Value = 'Denuncias'
$ws.Cells.Item(18, 2).Value = 'Denuncia'
$ws.Cells.Item(18, 3).Value = 'denuncia_otra_clasificacion'
$ws.Cells.Item(18, 4).Value = 'alfanumerico'
$ws.Cells.Item(18, 5).Value = 'Categoría o clasificación alternativa del hecho denunciado.'
$ws.Cells.Item(18, 6).Value = ""
$ws.Cells.Item(18, 7).Value = 'schema:Thing'

# Row 19: denuncia_fecha_ultimo_cambio
$ws.Cells.Item(19, 1).Value = 'Denuncias'
$ws.Cells.Item(19, 2).Value = 'Denuncia'
$ws.Cells.Item(19, 3).Value = 'denuncia_fecha_ultimo_cambio'
$ws.Cells.Item(19, 4).Value = 'fecha'
$ws.Cells.Item(19, 5).Value = 'Fecha de la última modificación que sufrió la denuncia.'
$ws.Cells.Item(19, 6).Value = "'" + '2019-03-15'
$ws.Cells.Item(19, 6).Style = "Normal"
$ws.Cells.Item(19, 7).Value = 'schema:Date'

# Row 20: denuncia_comentario_ultimo_cambio
$ws.Cells.Item(20, 1).Value = 'Denuncias'
$ws.Cells.Item(20, 2).Value = 'Denuncia'
$ws.Cells.Item(20, 3).Value = 'denuncia_comentario_ultimo_cambio'
$ws.Cells.Item(20, 4).Value = 'texto'
$ws.Cells.Item(20, 5).Value = 'Comentario agregado a la denuncia en el último cambio.'
$ws.Cells.Item(20, 6).Value = 'Se envió equipo a reparar el cable'
$ws.Cells.Item(20, 7).Value = 'schema:Thing'

# Remove now-obsolete rows 21-25 (their content was shifted up above).
$ws.Range("A21:H25").Delete()
